$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 6) by copying the formatting from row 5 first,
# so the new cells reuse the existing style indices (s="2"/s="3") instead
# of creating brand-new font/style entries.
$ws.Range("A5:L5").Copy()
$ws.Range("A6:L6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows(6).RowHeight = 18

# Populate the new row's values / formula.
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "{epoch_id}"
$ws.Range("C6").Value = "C00003"
$ws.Range("D6").Value = 444
$ws.Range("E6").Value = "6001-7890"
$ws.Range("F6").Value = "JPM4"
$ws.Range("G6").Value = "REDEMPTION"
$ws.Range("H6").Value = 46.53
$ws.Range("I6").Value = "USD"
$ws.Range("J6").Value = "{current_timestamp}"
$ws.Range("K6").Value = "BA1"
$ws.Range("L6").Formula = "=H6+50"

# Column D (bank.branch) is numeric-looking data that should be stored as
# text, so give it a text number format. Apply it across the whole column
# range (rows 2-6) which creates a single new text-formatted style reusing
# the existing bold-green font. (Must be done after all D-column values are
# set, otherwise re-assigning a value afterwards would coerce it to a text
# shared string instead of keeping the underlying numeric value.)
$ws.Range("D2:D6").NumberFormat = "@"

# Update the view: zoom to 93% and move the active selection to G6.
$ws.Application.ActiveWindow.Zoom = 93
$ws.Range("G6").Select()
